# Milestone 1, 30 oktober 2016
# Update the 21.10.2016 "Sitemap Wireframe" work-hours entry (E19) from 2 to
# 3.5 hours. E23 (=SUM(E3:E19)) recalculates automatically from 25.25 to
# 26.75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E19").Value = 3.5

# Move the active selection/cursor to reflect where the author left off
# editing (was E23, now E20).
$ws.Range("E20").Select()
